$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 441, shifting existing rows 441-462 down to 442-463.
$ws.Rows("441:441").Insert()

# Populate the newly inserted row 441 with the new record's data.
$ws.Range("A441").Value = 4
$ws.Range("B441").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C441").Value = "Los Lagos"
$ws.Range("D441").Value = 44939
$ws.Range("E441").Value = 10
$ws.Range("F441").Value = 100112023
$ws.Range("G441").Value = "Brócoli"
$ws.Range("H441").Value = "Sin especificar"
$ws.Range("I441").Value = "Primera"
$ws.Range("J441").Value = 1200
$ws.Range("K441").Value = 1500
$ws.Range("L441").Value = 1500
$ws.Range("M441").Value = 1500
$ws.Range("N441").Value = "$/unidad"
$ws.Range("O441").Value = "Región Metropolitana"
$ws.Range("P441").Value = 1500
$ws.Range("Q441").Value = 1
$ws.Range("R441").Value = "Hortaliza"
